$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header tweaks -----------------------------------------------
$ws.Range("F1").Value = "P_E ratio"
$ws.Range("G1").Value = "P_B ratio"

# Remove the old quarterly-results columns (K1:P1 headers, K2:P2 values)
# that used to sit to the right of the summary table; they're being
# relocated to a separate block starting at row 5.
$ws.Range("K1:P2").Clear()

# --- New quarterly-results block ---------------------------------------
# Header row (row 5)
$ws.Range("A5").Value = "Revenues"
$ws.Range("B5").Value = "Cost of Revenues"
$ws.Range("C5").Value = "General & Administrative Expenses in USD millions"
$ws.Range("D5").Value = "Operating Expenses in USD millions"
$ws.Range("E5").Value = "Interest Expense in USD millions"
$ws.Range("F5").Value = "Depreciation, Amortization & Accretion in USD millions"
$ws.Range("G5").Value = "EBITDA"
$ws.Range("H5").Value = "Gross Profit"
$ws.Range("I5").Value = "Net Income"
$ws.Range("J5").Value = "Weighted Average Shares"
$ws.Range("K5").Value = "Operating Income"

# Give row 5 the same bold/border/center style used by the row-1 header,
# by copying row 1's formatting over (reuses the existing style, rather
# than inventing a new one). A1 alone is used as the source (rather than
# A1:K1) because by this point K1 has already been cleared above.
$ws.Range("A1").Copy()
$ws.Range("A5:K5").PasteSpecial(-4122)  # xlPasteFormats

# Data row (row 6) - kept as plain text, matching how the rest of the
# sheet stores its numbers-with-commas.
$ws.Range("A6:K6").NumberFormat = "@"
$ws.Range("A6").Value = "11,953"
$ws.Range("B6").Value = "4,657"
$ws.Range("C6").Value = "3,667"
$ws.Range("D6").Value = "4,026"
$ws.Range("E6").Value = "368"
$ws.Range("F6").Value = "290"
$ws.Range("G6").Value = "4,199"
$ws.Range("H6").Value = "7,296"
$ws.Range("I6").Value = "3,087"
$ws.Range("J6").Value = "4,324"
$ws.Range("K6").Value = "3,270"

# Reset row 6's visual style back to the plain (unstyled) look used by
# row 2's data cells - the NumberFormat tweak above is just to force
# text storage, it shouldn't leave a lingering style behind. A2 alone is
# the source so every target cell (including the brand-new G6:K6 columns
# that don't line up with an existing column K..P) gets the same
# unstyled look.
$ws.Range("A2").Copy()
$ws.Range("A6:K6").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# --- Used range -----------------------------------------------------------
$ws.Range("A1").Select()
